# Auto-generated edit script: update cryptos list (prices / 1h volume %) per commit
# "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain-looking decimal strings as TEXT in the source file
# (t="inlineStr"). Pre-mark those cells as Text (NumberFormat "@") before writing the
# new value so Excel does not auto-coerce e.g. "357.11" into a Number, which would also
# silently drop formatting such as the trailing zero in "0.920" -> 0.92.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values (Price, Volume(1h), and the two re-ordered rows' Coin/Link)
$ws.Range("D2").Value = '51.962.17'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.778.80'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '357.11'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").Value = '109.44'
$ws.Range("E6").Value = '  -3.57%  '
$ws.Range("D7").Value = '0.564'
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").Value = '40.11'
$ws.Range("E10").Value = '  -3.53%  '
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '19.39'
$ws.Range("E13").Value = '  -3.19%  '
$ws.Range("D14").Value = '7.62'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '3.216.15'
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").Value = '2.795.28'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '0.931'
$ws.Range("E17").Value = '  +3.80%  '
$ws.Range("D18").Value = '51.833.80'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '3.15'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = '13.06'
$ws.Range("E21").Value = '  -3.85%  '
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '274.34'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").Value = '69.94'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '2.73'
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").Value = '26.56'
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '10.14'
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '0.143'
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.17'
$ws.Range("E30").Value = '  -3.44%  '
$ws.Range("D31").Value = '0.0466'
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").Value = '51.66'
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("D33").Value = '33.98'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("D35").Value = '0.0842'
$ws.Range("D36").Value = '5.26'
$ws.Range("E36").Value = '  +7.06%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '3.23'
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("D39").Value = '18.06'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '2.52'
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").Value = '2.23'
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("D44").Value = '121.05'
$ws.Range("E44").Value = '  -5.07%  '
$ws.Range("D45").Value = '21.88'
$ws.Range("E45").Value = '  -7.66%  '
$ws.Range("D46").Value = '2.068.07'
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("E47").Value = '  -3.08%  '
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").Value = '0.920'
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").Value = '8.96'
$ws.Range("E51").Value = '  +0.49%  '
